# Edit script: insert 5 new weekly price rows for
# "Comercializadora del Agro de Limarí - Tomate" ahead of the existing
# data block (old row 468 onward), shifting everything below down by 5
# rows (old A1:R551 -> new A1:R556), and fill the 5 new rows with the
# latest week's price data (fecha serial 44644).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank rows at the top of the existing data block; this pushes
# the old rows 468:551 down to 473:556 and keeps all of their values/
# formatting intact (including the date-formatted style on column D).
$ws.Rows("468:472").Insert()

# Row 468: Tomate / Larga vida / Primera
$ws.Range("A468").Value = 2
$ws.Range("B468").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C468").Value = "Coquimbo"
$ws.Range("D468").Value = 44644
$ws.Range("E468").Value = 4
$ws.Range("F468").Value = 100112020
$ws.Range("G468").Value = "Tomate"
$ws.Range("H468").Value = "Larga vida"
$ws.Range("I468").Value = "Primera"
$ws.Range("J468").Value = 800
$ws.Range("K468").Value = 9500
$ws.Range("L468").Value = 10000
$ws.Range("M468").Value = 9750
$ws.Range("N468").Value = "$/bandeja 18 kilos"
$ws.Range("O468").Value = "Provincia de Limarí"
$ws.Range("P468").Value = 542
$ws.Range("Q468").Value = 18
$ws.Range("R468").Value = "Hortaliza"

# Row 469: Tomate / Larga vida / Segunda
$ws.Range("A469").Value = 2
$ws.Range("B469").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C469").Value = "Coquimbo"
$ws.Range("D469").Value = 44644
$ws.Range("E469").Value = 4
$ws.Range("F469").Value = 100112020
$ws.Range("G469").Value = "Tomate"
$ws.Range("H469").Value = "Larga vida"
$ws.Range("I469").Value = "Segunda"
$ws.Range("J469").Value = 500
$ws.Range("K469").Value = 7500
$ws.Range("L469").Value = 8000
$ws.Range("M469").Value = 7750
$ws.Range("N469").Value = "$/bandeja 18 kilos"
$ws.Range("O469").Value = "Provincia de Limarí"
$ws.Range("P469").Value = 431
$ws.Range("Q469").Value = 18
$ws.Range("R469").Value = "Hortaliza"

# Row 470: Tomate / Larga vida / Tercera
$ws.Range("A470").Value = 2
$ws.Range("B470").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C470").Value = "Coquimbo"
$ws.Range("D470").Value = 44644
$ws.Range("E470").Value = 4
$ws.Range("F470").Value = 100112020
$ws.Range("G470").Value = "Tomate"
$ws.Range("H470").Value = "Larga vida"
$ws.Range("I470").Value = "Tercera"
$ws.Range("J470").Value = 400
$ws.Range("K470").Value = 5500
$ws.Range("L470").Value = 6000
$ws.Range("M470").Value = 5750
$ws.Range("N470").Value = "$/bandeja 18 kilos"
$ws.Range("O470").Value = "Provincia de Limarí"
$ws.Range("P470").Value = 319
$ws.Range("Q470").Value = 18
$ws.Range("R470").Value = "Hortaliza"

# Row 471: Tomate / Semiduro / Primera
$ws.Range("A471").Value = 2
$ws.Range("B471").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C471").Value = "Coquimbo"
$ws.Range("D471").Value = 44644
$ws.Range("E471").Value = 4
$ws.Range("F471").Value = 100112020
$ws.Range("G471").Value = "Tomate"
$ws.Range("H471").Value = "Semiduro"
$ws.Range("I471").Value = "Primera"
$ws.Range("J471").Value = 500
$ws.Range("K471").Value = 7500
$ws.Range("L471").Value = 8000
$ws.Range("M471").Value = 7750
$ws.Range("N471").Value = "$/bandeja 18 kilos"
$ws.Range("O471").Value = "Provincia de Limarí"
$ws.Range("P471").Value = 431
$ws.Range("Q471").Value = 18
$ws.Range("R471").Value = "Hortaliza"

# Row 472: Tomate / Semiduro / Segunda
$ws.Range("A472").Value = 2
$ws.Range("B472").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C472").Value = "Coquimbo"
$ws.Range("D472").Value = 44644
$ws.Range("E472").Value = 4
$ws.Range("F472").Value = 100112020
$ws.Range("G472").Value = "Tomate"
$ws.Range("H472").Value = "Semiduro"
$ws.Range("I472").Value = "Segunda"
$ws.Range("J472").Value = 300
$ws.Range("K472").Value = 5500
$ws.Range("L472").Value = 6000
$ws.Range("M472").Value = 5750
$ws.Range("N472").Value = "$/bandeja 18 kilos"
$ws.Range("O472").Value = "Provincia de Limarí"
$ws.Range("P472").Value = 319
$ws.Range("Q472").Value = 18
$ws.Range("R472").Value = "Hortaliza"
